# Depot log v1.7 IKAMR 2020-03-30 — part 2
#
# The commit normalises three "Lastet ned ...zip/.tar" / "Pakket ut ...tar"
# log-entry sentences so that every placeholder filename reads
# "<SIP uuid>" instead of the inconsistent "abcdef" / "<uuid>" forms.

$d = $word.ActiveDocument

# 1) "Lastet ned abcdef.zip ..."  ->  "Lastet ned <SIP uuid>.zip ..."
$d.Content.Find.Execute(
    "Lastet ned abcdef.zip på n.n GB fra Produsent",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lastet ned <SIP uuid>.zip på n.n GB fra Produsent", 2
)

# 2) "Pakket ut <uuid>.tar ..."  ->  "Pakket ut <SIP uuid>.tar ..."
$d.Content.Find.Execute(
    "Pakket ut <uuid>.tar på n.n GB",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Pakket ut <SIP uuid>.tar på n.n GB", 2
)

# 3) "Lastet ned <uuid>.tar ..."  ->  "Lastet ned <SIP uuid>.tar ..."
$d.Content.Find.Execute(
    "Lastet ned <uuid>.tar på n.n GB fra Produsent",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Lastet ned <SIP uuid>.tar på n.n GB fra Produsent", 2
)
